# Update "想去人数" (wish-to-go count) figures on both the "展览" and
# "全部类型" sheets to reflect newly scraped totals.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 895
    $ws.Range("F3").Value = 4533
    $ws.Range("F5").Value = 792
}
